$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 (id=45)
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "RemoteJob"
$ws.Cells.Item(46, 3).Value = "2017-11-12 09:03:38"
$ws.Cells.Item(46, 4).Value = "2017-11-12 09:03:38"
$ws.Cells.Item(46, 5).Value = 569
$ws.Cells.Item(46, 6).Value = 473
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = 1
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0.33
$ws.Cells.Item(46, 11).Value = 36
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 0

# Row 47 (id=46)
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "npm start"
$ws.Cells.Item(47, 3).Value = "2017-11-12 09:03:44"
$ws.Cells.Item(47, 4).Value = "2017-11-12 09:03:48"
$ws.Cells.Item(47, 5).Value = 420
$ws.Cells.Item(47, 6).Value = 416
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 36
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0

# Row 48 (id=47)
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "npm start"
$ws.Cells.Item(48, 3).Value = "2017-11-12 09:27:23"
$ws.Cells.Item(48, 4).Value = "2017-11-12 09:27:30"
$ws.Cells.Item(48, 5).Value = 537
$ws.Cells.Item(48, 6).Value = 494
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 37
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 0

# Row 49 (id=48)
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "RemoteJob"
$ws.Cells.Item(49, 3).Value = "2017-11-12 09:27:30"
$ws.Cells.Item(49, 4).Value = "2017-11-12 09:27:33"
$ws.Cells.Item(49, 5).Value = 549
$ws.Cells.Item(49, 6).Value = 472
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = 1
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 38
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = 0

# Row 50 (id=49)
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "npm start"
$ws.Cells.Item(50, 3).Value = "2017-11-12 09:27:52"
$ws.Cells.Item(50, 4).Value = "2017-11-12 09:27:53"
$ws.Cells.Item(50, 5).Value = 413
$ws.Cells.Item(50, 6).Value = 358
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 38
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0

# Row 51 (id=50)
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "RemoteJob"
$ws.Cells.Item(51, 3).Value = "2017-11-12 09:27:53"
$ws.Cells.Item(51, 4).Value = "2017-11-12 09:28:00"
$ws.Cells.Item(51, 5).Value = 541
$ws.Cells.Item(51, 6).Value = 384
$ws.Cells.Item(51, 7).Value = 7
$ws.Cells.Item(51, 8).Value = 1
$ws.Cells.Item(51, 9).Value = 12
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 39
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0

# Row 52 (id=51)
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "npm start"
$ws.Cells.Item(52, 3).Value = "2017-11-12 09:28:37"
$ws.Cells.Item(52, 4).Value = "2017-11-12 09:28:39"
$ws.Cells.Item(52, 5).Value = 589
$ws.Cells.Item(52, 6).Value = 144
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 11).Value = 39
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 13).Value = 0

# Row 53 (id=52)
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "RemoteJob"
$ws.Cells.Item(53, 3).Value = "2017-11-12 09:28:39"
$ws.Cells.Item(53, 4).Value = "2017-11-12 09:28:51"
$ws.Cells.Item(53, 5).Value = 566
$ws.Cells.Item(53, 6).Value = 313
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 1
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 40
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0

# Row 54 (id=53)
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "npm start"
$ws.Cells.Item(54, 3).Value = "2017-11-12 09:29:12"
$ws.Cells.Item(54, 4).Value = "2017-11-12 09:29:14"
$ws.Cells.Item(54, 5).Value = 615
$ws.Cells.Item(54, 6).Value = 328
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 40
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = 0

# Row 55 (id=54)
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "RemoteJob"
$ws.Cells.Item(55, 3).Value = "2017-11-12 09:29:14"
$ws.Cells.Item(55, 4).Value = "2017-11-12 09:29:18"
$ws.Cells.Item(55, 5).Value = 593
$ws.Cells.Item(55, 6).Value = 371
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 41
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = 0

# Row 56 (id=55)
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "npm start"
$ws.Cells.Item(56, 3).Value = "2017-11-12 09:29:37"
$ws.Cells.Item(56, 4).Value = "2017-11-12 09:29:38"
$ws.Cells.Item(56, 5).Value = 347
$ws.Cells.Item(56, 6).Value = 388
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 1
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0.33
$ws.Cells.Item(56, 11).Value = 41
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = 0

# Row 57 (id=56)
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "RemoteJob"
$ws.Cells.Item(57, 3).Value = "2017-11-12 09:29:38"
$ws.Cells.Item(57, 4).Value = "2017-11-12 09:30:07"
$ws.Cells.Item(57, 5).Value = 373
$ws.Cells.Item(57, 6).Value = 278
$ws.Cells.Item(57, 7).Value = 8
$ws.Cells.Item(57, 8).Value = 1
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0.33
$ws.Cells.Item(57, 11).Value = 42
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = 0

# Row 58 (id=57)
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "npm start"
$ws.Cells.Item(58, 3).Value = "2017-11-12 09:30:44"
$ws.Cells.Item(58, 4).Value = "2017-11-12 09:30:46"
$ws.Cells.Item(58, 5).Value = 510
$ws.Cells.Item(58, 6).Value = 393
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 42
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = 0

# Row 59 (id=58)
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = "npm start"
$ws.Cells.Item(59, 3).Value = "2017-11-12 09:31:11"
$ws.Cells.Item(59, 4).Value = "2017-11-12 09:31:13"
$ws.Cells.Item(59, 5).Value = 503
$ws.Cells.Item(59, 6).Value = 357
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 11).Value = 43
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 13).Value = 0

# Row 60 (id=59)
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "npm start"
$ws.Cells.Item(60, 3).Value = "2017-11-12 09:37:17"
$ws.Cells.Item(60, 4).Value = "2017-11-12 09:37:18"
$ws.Cells.Item(60, 5).Value = 518
$ws.Cells.Item(60, 6).Value = 384
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 44
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = 0

# Row 61 (id=60)
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "npm start"
$ws.Cells.Item(61, 3).Value = "2017-11-12 09:37:40"
$ws.Cells.Item(61, 4).Value = "2017-11-12 09:37:42"
$ws.Cells.Item(61, 5).Value = 547
$ws.Cells.Item(61, 6).Value = 505
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 45
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = 0

# Row 62 (id=61)
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "npm start"
$ws.Cells.Item(62, 3).Value = "2017-11-12 09:38:50"
$ws.Cells.Item(62, 4).Value = "2017-11-12 09:38:52"
$ws.Cells.Item(62, 5).Value = 498
$ws.Cells.Item(62, 6).Value = 413
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 46
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = 0

# Row 63 (id=62)
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "npm start"
$ws.Cells.Item(63, 3).Value = "2017-11-12 09:39:55"
$ws.Cells.Item(63, 4).Value = "2017-11-12 09:39:59"
$ws.Cells.Item(63, 5).Value = 474
$ws.Cells.Item(63, 6).Value = 419
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 47
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0

# Row 64 (id=63)
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "npm start"
$ws.Cells.Item(64, 3).Value = "2017-11-12 09:40:13"
$ws.Cells.Item(64, 4).Value = "2017-11-12 09:40:16"
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 48
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0

# Row 65 (id=64)
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "npm start"
$ws.Cells.Item(65, 3).Value = "2017-11-12 09:40:31"
$ws.Cells.Item(65, 4).Value = "2017-11-12 09:40:34"
$ws.Cells.Item(65, 5).Value = 483
$ws.Cells.Item(65, 6).Value = 469
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 49
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0

# Row 66 (id=65)
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "npm start"
$ws.Cells.Item(66, 3).Value = "2017-11-12 09:40:44"
$ws.Cells.Item(66, 4).Value = "2017-11-12 09:40:53"
$ws.Cells.Item(66, 5).Value = 629
$ws.Cells.Item(66, 6).Value = 206
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 50
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = 0

# Row 67 (id=66)
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "RemoteJob"
$ws.Cells.Item(67, 3).Value = "2017-11-12 09:40:53"
$ws.Cells.Item(67, 4).Value = "2017-11-12 09:41:21"
$ws.Cells.Item(67, 5).Value = 465
$ws.Cells.Item(67, 6).Value = 93
$ws.Cells.Item(67, 7).Value = 4
$ws.Cells.Item(67, 8).Value = 1
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 1.32
$ws.Cells.Item(67, 11).Value = 51
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = 0

# Row 68 (id=67)
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "book-manager"
$ws.Cells.Item(68, 3).Value = "2017-11-12 09:41:21"
$ws.Cells.Item(68, 4).Value = "2017-11-12 09:41:26"
$ws.Cells.Item(68, 5).Value = 1059
$ws.Cells.Item(68, 6).Value = 10
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 1
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0.33
$ws.Cells.Item(68, 11).Value = 51
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0

# Row 69 (id=68)
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "RemoteJob"
$ws.Cells.Item(69, 3).Value = "2017-11-12 09:41:26"
$ws.Cells.Item(69, 4).Value = "2017-11-12 09:41:29"
$ws.Cells.Item(69, 5).Value = 469
$ws.Cells.Item(69, 6).Value = 89
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 1
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0.33
$ws.Cells.Item(69, 11).Value = 51
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0

# Row 70 (id=69)
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "book-manager"
$ws.Cells.Item(70, 3).Value = "2017-11-12 09:41:29"
$ws.Cells.Item(70, 4).Value = "2017-11-12 09:41:30"
$ws.Cells.Item(70, 5).Value = 1061
$ws.Cells.Item(70, 6).Value = 12
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 1
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0.33
$ws.Cells.Item(70, 11).Value = 51
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0

# Row 71 (id=70)
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "RemoteJob"
$ws.Cells.Item(71, 3).Value = "2017-11-12 09:41:30"
$ws.Cells.Item(71, 4).Value = "2017-11-12 09:41:32"
$ws.Cells.Item(71, 5).Value = 379
$ws.Cells.Item(71, 6).Value = 389
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 1
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0.33
$ws.Cells.Item(71, 11).Value = 51
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0

# Row 72 (id=71)
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "npm start"
$ws.Cells.Item(72, 3).Value = "2017-11-12 09:41:52"
$ws.Cells.Item(72, 4).Value = "2017-11-12 09:41:57"
$ws.Cells.Item(72, 5).Value = 397
$ws.Cells.Item(72, 6).Value = 546
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 51
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = 0

# Row 73 (id=72)
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = "npm start"
$ws.Cells.Item(73, 3).Value = "2017-11-12 09:46:47"
$ws.Cells.Item(73, 4).Value = "2017-11-12 09:46:50"
$ws.Cells.Item(73, 5).Value = 408
$ws.Cells.Item(73, 6).Value = 110
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 52
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0

# Row 74 (id=73)
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = "book-manager"
$ws.Cells.Item(74, 3).Value = "2017-11-12 10:06:41"
$ws.Cells.Item(74, 4).Value = "2017-11-12 10:06:43"
$ws.Cells.Item(74, 5).Value = 351
$ws.Cells.Item(74, 6).Value = 173
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 53
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
